$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tbl = $shape.Table
$cell = $tbl.Cell(1, 2)
$tr = $cell.Shape.TextFrame.TextRange
$tr.Text = "retro Genes"
